$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'303.13"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'5.59%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'31.99"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'9.61%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.250"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'1.04%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07459"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'6.90%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'7.849"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'5.55%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'3.796"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'7.01%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'1.530"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'8.92%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.9186"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'1.96%"
$ws.Range("E9").Style = "Normal"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1680"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'4.27%"
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").Value = "'0.08040"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'8.31%"
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.07944"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'2.28%"
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03031"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'3.14%"
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09920"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'9.93%"
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001492"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-6.17%"
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = "CoinExToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D16").Value = "'0.04618"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'2.23%"
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D17").Value = "'0.006322"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'1.61%"
$ws.Range("E17").Style = "Normal"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D18").Value = "'3.462"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-0.30%"
$ws.Range("E18").Style = "Normal"
$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D19").Value = "'2.234"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'0.21%"
$ws.Range("E19").Style = "Normal"
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20").Value = "'0.3326"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'2.65%"
$ws.Range("E20").Style = "Normal"
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").Value = "'0.1338"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'0.32%"
$ws.Range("E21").Style = "Normal"
$ws.Range("B22").Value = "MCDex"
$ws.Range("C22").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D22").Value = "'4.476"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'11.01%"
$ws.Range("E22").Style = "Normal"
$ws.Range("B23").Value = "ZBToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D23").Value = "'0.1621"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'1.42%"
$ws.Range("E23").Style = "Normal"
$ws.Range("B24").Value = "BitKan"
$ws.Range("C24").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D24").Value = "'0.001217"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.70%"
$ws.Range("E24").Style = "Normal"
$ws.Range("B25").Value = "HotbitToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D25").Value = "'0.004446"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'7.43%"
$ws.Range("E25").Style = "Normal"
$ws.Range("B26").Value = "NitroEx"
$ws.Range("C26").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D26").Value = "'0.0001399"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'19.72%"
$ws.Range("E26").Style = "Normal"
$ws.Range("B27").Value = "UpBots"
$ws.Range("C27").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D27").Value = "'0.0001749"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'4.88%"
$ws.Range("E27").Style = "Normal"
$ws.Range("B28").Value = "Spectre.aiUtilityToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/yxQ8LoZvwJ6Ee+spectreaiutilitytoken-sxut"
$ws.Range("D28").Value = "'--"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'--%"
$ws.Range("E28").Style = "Normal"
$ws.Range("B29").Value = "LegolasExchange"
$ws.Range("C29").Value = "https://coinranking.com/coin/zEMEnlPs_94tc+legolasexchange-lgo"
$ws.Range("B30").Value = "BitZToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/nLHHPBKDJSEee+bitztoken-bz"
$ws.Range("B31").Value = "Birake"
$ws.Range("C31").Value = "https://coinranking.com/coin/dTOfofFqKQiY5+birake-bir"
$ws.Range("B32").Value = "NashExchange"
$ws.Range("C32").Value = "https://coinranking.com/coin/9LcSTo0q-+nashexchange-nex"
$ws.Range("B33").Value = "AAXToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/LNePqkIhk+aaxtoken-aab"
$ws.Range("B34").Value = "CenX"
$ws.Range("C34").Value = "https://coinranking.com/coin/V4XJUvLQb+cenx-cenx"
$ws.Range("B35").Value = "BNIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/n194X9uHp+bnixtoken-bnix"
$ws.Range("B36").Value = "Polkally"
$ws.Range("C36").Value = "https://coinranking.com/coin/NkDWUL8F-+polkally-kally"
$ws.Range("B37").Value = "Charli3"
$ws.Range("C37").Value = "https://coinranking.com/coin/8SgjMSqUk+charli3-c3"
$ws.Range("B38").Value = "BlubitexToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Y9oImHIW5+blubitextoken-bbe"
$ws.Range("B39").Value = "One"
$ws.Range("C39").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D39").Value = "'0.01744"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'2,581.14%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.04489"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'2.86%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007180"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'3.63%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1347"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'8.02%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002199"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'6.33%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.01283"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'10.77%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00006135"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'4.92%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.7097"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-63.21%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.01301"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-0.41%"
$ws.Range("E47").Style = "Normal"
